$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "2024.08.01. 월간"
$ws.Range("B7").Value = "2024년 09월 04일 19시 07분 39초"
$ws.Range("C10").Value = "121"
$ws.Range("D10").Value = "29.95"
$ws.Range("C11").Value = "283"
$ws.Range("D11").Value = "70.05"
$ws.Range("C14").Value = "8"
$ws.Range("D14").Value = "1.98"
$ws.Range("C15").Value = "30"
$ws.Range("D15").Value = "7.43"
$ws.Range("C16").Value = "56"
$ws.Range("D16").Value = "13.86"
$ws.Range("C17").Value = "111"
$ws.Range("D17").Value = "27.48"
$ws.Range("C18").Value = "20"
$ws.Range("D18").Value = "4.95"
$ws.Range("C19").Value = "44"
$ws.Range("D19").Value = "10.89"
$ws.Range("C20").Value = "9"
$ws.Range("D20").Value = "2.23"
$ws.Range("C21").Value = "22"
$ws.Range("D21").Value = "5.45"
$ws.Range("D22").Value = "1.49"
$ws.Range("C23").Value = "14"
$ws.Range("D23").Value = "3.47"
$ws.Range("C24").Value = "7"
$ws.Range("D24").Value = "1.73"
$ws.Range("C25").Value = "19"
$ws.Range("D25").Value = "4.7"
$ws.Range("C26").Value = "2"
$ws.Range("D26").Value = "0.5"
$ws.Range("C27").Value = "15"
$ws.Range("D27").Value = "3.71"
$ws.Range("C28").Value = "4"
$ws.Range("D28").Value = "0.99"
$ws.Range("C29").Value = "7"
$ws.Range("D29").Value = "1.73"
$ws.Range("C30").Value = "5"
$ws.Range("D30").Value = "1.24"
$ws.Range("C31").Value = "9"
$ws.Range("D31").Value = "2.23"
$ws.Range("C32").Value = "4"
$ws.Range("D32").Value = "0.99"
$ws.Range("C33").Value = "12"
$ws.Range("D33").Value = "2.97"
